$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53: No Accounting for Waste | Enchanted Electrum Ink
$ws.Range("H53").Value = 224
$ws.Range("I53").Value = 210.85715
$ws.Range("K53").Value = 210.85715
$ws.Range("M53").Value = 426.14285

# Row 86: Filling in the Blanks | Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 5036.8604
$ws.Range("I86").Value = 4481.607
$ws.Range("J86").Value = 6073.3335
$ws.Range("K86").Value = 4481.607
$ws.Range("L86").Value = 6073.3335
$ws.Range("M86").Value = -3358.607
$ws.Range("N86").Value = -8319.333500000001

# Row 89: Ink into Antiquity (L) | Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 5036.8604
$ws.Range("I89").Value = 4481.607
$ws.Range("J89").Value = 6073.3335
$ws.Range("K89").Value = 22408.035
$ws.Range("L89").Value = 30366.6675
$ws.Range("M89").Value = -16792.035
$ws.Range("N89").Value = -41598.6675

# Row 92: Whinier than the Sword | Enchanted Koppranickel Ink
$ws.Range("H92").Value = 1124.6957
$ws.Range("I92").Value = 1245.7894
$ws.Range("K92").Value = 1245.7894
$ws.Range("M92").Value = 2.210600000000113

# Row 98: The Dotted Line | Enchanted Durium Ink
$ws.Range("H98").Value = 2227.28
$ws.Range("I98").Value = 2190
$ws.Range("J98").Value = 2423
$ws.Range("K98").Value = 2190
$ws.Range("L98").Value = 2423
$ws.Range("M98").Value = -692
$ws.Range("N98").Value = -5419

# Row 99: Rumor Has It | Commanding Craftsman's Tea
$ws.Range("H99").Value = 539
$ws.Range("I99").Value = 539
$ws.Range("K99").Value = 1617
$ws.Range("M99").Value = -119

# Row 100: Asking for a Friend | Beetle Glue
$ws.Range("H100").Value = 3432.1365
$ws.Range("J100").Value = 5054.364
$ws.Range("L100").Value = 5054.364
$ws.Range("N100").Value = -6136.364

# Row 122: Wishful Inking | Enchanted High Durium Ink
$ws.Range("H122").Value = 2227.28
$ws.Range("I122").Value = 2190
$ws.Range("J122").Value = 2423
$ws.Range("K122").Value = 6570
$ws.Range("L122").Value = 7269
$ws.Range("M122").Value = -4120
$ws.Range("N122").Value = -12169

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff | Cobalt Ingot
$ws.Range("H61").Value = 8336355.5
$ws.Range("I61").Value = 2696.4443
$ws.Range("K61").Value = 2696.4443
$ws.Range("M61").Value = -2484.4443

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 741127.3
$ws.Range("I74").Value = 838321.1
$ws.Range("J74").Value = 19115.857
$ws.Range("K74").Value = 838321.1
$ws.Range("L74").Value = 19115.857
$ws.Range("M74").Value = -837447.1
$ws.Range("N74").Value = -20863.857

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 741127.3
$ws.Range("I77").Value = 838321.1
$ws.Range("J77").Value = 19115.857
$ws.Range("K77").Value = 4191605.5
$ws.Range("L77").Value = 95579.285
$ws.Range("M77").Value = -4187237.5
$ws.Range("N77").Value = -104315.285

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Range("H132").Value = 4728.0938
$ws.Range("I132").Value = 2910.2
$ws.Range("J132").Value = 5554.409
$ws.Range("K132").Value = 8730.599999999999
$ws.Range("L132").Value = 16663.227
$ws.Range("M132").Value = -6200.599999999999
$ws.Range("N132").Value = -21723.227

# Row 136: Metal with Mettle | Cobalt Tungsten Ingot
$ws.Range("H136").Value = 8336355.5
$ws.Range("I136").Value = 2696.4443
$ws.Range("K136").Value = 8089.3329
$ws.Range("M136").Value = -5539.3329

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt | Iron Ingot
$ws.Range("H20").Value = 41695.023
$ws.Range("I20").Value = 48015.184
$ws.Range("J20").Value = 20838.5
$ws.Range("K20").Value = 48015.184
$ws.Range("L20").Value = 20838.5
$ws.Range("M20").Value = -47768.184
$ws.Range("N20").Value = -21332.5

# Row 58: You Stay on That Side | Cobalt Pliers
$ws.Range("H58").Value = 71340
$ws.Range("J58").Value = 71340
$ws.Range("L58").Value = 71340
$ws.Range("N58").Value = -71928

# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 5701.6665
$ws.Range("I86").Value = 3548.5
$ws.Range("J86").Value = 6778.25
$ws.Range("K86").Value = 3548.5
$ws.Range("L86").Value = 6778.25
$ws.Range("M86").Value = -2425.5
$ws.Range("N86").Value = -9024.25

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 5701.6665
$ws.Range("I89").Value = 3548.5
$ws.Range("J89").Value = 6778.25
$ws.Range("K89").Value = 17742.5
$ws.Range("L89").Value = 33891.25
$ws.Range("M89").Value = -12126.5
$ws.Range("N89").Value = -45123.25

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent | Maple Lumber
$ws.Range("H7").Value = 86.86957
$ws.Range("I7").Value = 49.42857
$ws.Range("J7").Value = 145.11111
$ws.Range("K7").Value = 49.42857
$ws.Range("L7").Value = 145.11111
$ws.Range("M7").Value = 63.57143
$ws.Range("N7").Value = -371.11111

# Row 16: Raise the Roof | Ash Lumber
$ws.Range("H16").Value = 3701
$ws.Range("I16").Value = 3685.1428
$ws.Range("K16").Value = 3685.1428
$ws.Range("M16").Value = -3398.1428

# Row 31: Wall Not Found | Walnut Lumber
$ws.Range("H31").Value = 7781166.5
$ws.Range("I31").Value = 5001749.5
$ws.Range("J31").Value = 10004700
$ws.Range("K31").Value = 5001749.5
$ws.Range("L31").Value = 10004700
$ws.Range("M31").Value = -5001454.5
$ws.Range("N31").Value = -10005290

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Range("H34").Value = 7781166.5
$ws.Range("I34").Value = 5001749.5
$ws.Range("J34").Value = 10004700
$ws.Range("K34").Value = 5001749.5
$ws.Range("L34").Value = 10004700
$ws.Range("M34").Value = -5001547.5
$ws.Range("N34").Value = -10005104

# Row 86: Birch, Please | Birch Lumber
$ws.Range("H86").Value = 36428.97
$ws.Range("I86").Value = 103552.75
$ws.Range("J86").Value = 13081.565
$ws.Range("K86").Value = 103552.75
$ws.Range("L86").Value = 13081.565
$ws.Range("M86").Value = -102429.75
$ws.Range("N86").Value = -15327.565

# Row 89: Built This City on Blocks and Soul (L) | Birch Lumber
$ws.Range("H89").Value = 36428.97
$ws.Range("I89").Value = 103552.75
$ws.Range("J89").Value = 13081.565
$ws.Range("K89").Value = 517763.75
$ws.Range("L89").Value = 65407.825
$ws.Range("M89").Value = -512147.75
$ws.Range("N89").Value = -76639.82500000001

# Row 113: Patient Patients | White Ash Lumber
$ws.Range("H113").Value = 3701
$ws.Range("I113").Value = 3685.1428
$ws.Range("K113").Value = 3685.1428
$ws.Range("M113").Value = -1515.1428

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 2210.7856
$ws.Range("I134").Value = 2363.5144
$ws.Range("K134").Value = 7090.5432
$ws.Range("M134").Value = -4555.5432

$ws = $wb.Worksheets.Item("CUL")
# Row 2: Pork Is a Salty Food | Table Salt
$ws.Range("H2").Value = 78.64
$ws.Range("I2").Value = 16.294117
$ws.Range("J2").Value = 211.125
$ws.Range("K2").Value = 97.764702
$ws.Range("L2").Value = 1266.75
$ws.Range("M2").Value = 15.235298
$ws.Range("N2").Value = -1492.75

# Row 7: It's Always Sunny in Vylbrand | Raisins
$ws.Range("H7").Value = 247.6
$ws.Range("J7").Value = 400
$ws.Range("L7").Value = 1200
$ws.Range("N7").Value = -1424

# Row 38: Pretty as a Picture | Dark Vinegar
$ws.Range("H38").Value = 152.1579
$ws.Range("I38").Value = 253.55556
$ws.Range("K38").Value = 760.66668
$ws.Range("M38").Value = -413.66668

# Row 63: The Next to Last Supper | Stuffed Cabbage Rolls
$ws.Range("H63").Value = 2777.25
$ws.Range("I63").Value = 2943.6
$ws.Range("K63").Value = 8830.799999999999
$ws.Range("M63").Value = -8081.799999999999

# Row 64: The Aroma of Faith | Baked Onion Soup
$ws.Range("H64").Value = 6946.3
$ws.Range("I64").Value = 3928.4285
$ws.Range("K64").Value = 11785.2855
$ws.Range("M64").Value = -11515.2855

# Row 66: Nostalgia through the Stomach (L) | Stuffed Cabbage Rolls
$ws.Range("H66").Value = 2777.25
$ws.Range("I66").Value = 2943.6
$ws.Range("K66").Value = 26492.4
$ws.Range("M66").Value = -22748.4

# Row 67: Soup's On (L) | Baked Onion Soup
$ws.Range("H67").Value = 6946.3
$ws.Range("I67").Value = 3928.4285
$ws.Range("K67").Value = 11785.2855
$ws.Range("M67").Value = -10849.2855

# Row 87: Soup That Eats Like a Knight | Clam Chowder
$ws.Range("H87").Value = 17667.166
$ws.Range("I87").Value = 5751.5
$ws.Range("K87").Value = 17254.5
$ws.Range("M87").Value = -16006.5

# Row 90: Like Ma Used to Make (L) | Clam Chowder
$ws.Range("H90").Value = 17667.166
$ws.Range("I90").Value = 5751.5
$ws.Range("K90").Value = 51763.5
$ws.Range("M90").Value = -45523.5

# Row 114: One Last Meal | Mushroom Saute
$ws.Range("H114").Value = 1760.1111
$ws.Range("I114").Value = 210.75
$ws.Range("K114").Value = 632.25
$ws.Range("M114").Value = 2621.75

# Row 134: Don't Knock It Till You've Tried It | Mezcal-marinated Swampmonk
$ws.Range("H134").Value = 5783.857
$ws.Range("I134").Value = 2230.7334
$ws.Range("K134").Value = 6692.2002
$ws.Range("M134").Value = -1622.2002

# Row 139: Najoothie | Wild Banana Blend
$ws.Range("H139").Value = 7317.76
$ws.Range("I139").Value = 4309
$ws.Range("J139").Value = 12666.667
$ws.Range("K139").Value = 12927
$ws.Range("L139").Value = 38000.001
$ws.Range("M139").Value = -7787
$ws.Range("N139").Value = -48280.001

# Row 140: Sweet, Sweet Bean Juice | Mesquite Juice
$ws.Range("H140").Value = 3439.2727
$ws.Range("I140").Value = 2883.2
$ws.Range("K140").Value = 8649.599999999999
$ws.Range("M140").Value = -3469.599999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit | Mythrite Ingot
$ws.Range("H70").Value = 15476.2
$ws.Range("I70").Value = 24291.3
$ws.Range("J70").Value = 9599.467000000001
$ws.Range("K70").Value = 24291.3
$ws.Range("L70").Value = 9599.467000000001
$ws.Range("M70").Value = -24021.3
$ws.Range("N70").Value = -10139.467

# Row 73: Hulls of Broken Dreams (L) | Mythrite Ingot
$ws.Range("H73").Value = 15476.2
$ws.Range("I73").Value = 24291.3
$ws.Range("J73").Value = 9599.467000000001
$ws.Range("K73").Value = 24291.3
$ws.Range("L73").Value = 9599.467000000001
$ws.Range("M73").Value = -23355.3
$ws.Range("N73").Value = -11471.467

# Row 107: Whetstones for the Workers | Hard Mudstone Whetstone
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1000
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 920
$ws.Range("N107").ClearContents()

# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 2817.6538
$ws.Range("I122").Value = 2978.6
$ws.Range("J122").Value = 2281.1667
$ws.Range("K122").Value = 8935.799999999999
$ws.Range("L122").Value = 6843.500100000001
$ws.Range("M122").Value = -6485.799999999999
$ws.Range("N122").Value = -11743.5001

$ws = $wb.Worksheets.Item("WVR")
# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 11776948
$ws.Range("I136").Value = 5436566.5
$ws.Range("K136").Value = 16309699.5
$ws.Range("M136").Value = -16307149.5

